$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New falling-sensor dataset (rows 2-31): the commit replaces the old 20-row
# block with 13 freshly-sampled rows, then re-appends the previous rows 2-18
# (timestamps recomputed) so the sheet grows from 20 to 30 data rows.
$data = @(
    ,(0, "falling", -0.7822842597961426, 0.7009620666503906, 0.2034818530082702, -0.02587450614997311, 0.008312130346894184, -0.03830997752291809)
    ,(100, "falling", -0.6724414825439453, 0.7283806800842285, 0.0338471233844757, -0.03166318188110984, 0.005548692618807229, 0.004886921495199072)
    ,(200, "falling", -0.955592155456543, 0.6646170616149902, 0.1322764456272125, -0.01120646846746765, -0.01967131506119446, 0.04684027514996975)
    ,(300, "falling", -0.5998220443725586, 0.6136178970336914, -0.3951224088668823, -0.005999569415247833, -0.01865320652723307, 0.03610649971025327)
    ,(400, "falling", -0.6383709907531738, 0.829352855682373, 0.2481068670749664, -0.01938042674391034, 0.02639810448246341, -0.02997603196473341)
    ,(500, "falling", -0.8866777420043945, 0.6802849769592285, 0.1837458908557891, 0.04876013357369639, 0.03156137033476548, -0.02466732121649261)
    ,(600, "falling", -0.7689142227172852, 0.6263480186462402, -0.1205646991729736, 0.07661995095466947, 0.006654067996091003, 0.02539454134447233)
    ,(700, "falling", -0.7416219711303711, 0.6671648025512695, -0.0121779441833496, -0.008406669134274104, 0.01747510954737656, 0.02251474718962388)
    ,(800, "falling", -0.9570107460021972, 0.7521457672119141, 0.0322656333446502, -0.01705332100391391, 0.02540181328852965, 0.02260201397196701)
    ,(900, "falling", -0.8430213928222656, 0.7305541038513184, 0.0520271658897399, 0.004734205614243124, 0.03287764106478005, -0.002007128720703896)
    ,(1000, "falling", -0.7825779914855957, 0.6732511520385742, -0.1082392036914825, 0.002516182849094959, 0.03615013385812434, -0.01502437632353524)
    ,(1100, "falling", -0.6774282455444336, 0.5100307464599609, -0.1018824279308319, -0.01514073128678966, 0.0402807449655873, -0.0335248665263255)
    ,(1200, "falling", -0.801264762878418, 0.6673240661621094, 0.1666318476200103, 0.004646939025925644, 0.05990842623370023, -0.04932736818279529)
    ,(1300, "falling", -0.7477993965148926, 0.8460922241210938, 0.0258597135543823, 0.01418807215633853, 0.06712245657330498, -0.07486735071454727)
    ,(1400, "falling", -0.8064756393432617, 0.8524413108825684, -0.0703473389148712, -0.02585268907603775, -0.07008951618557867, -0.06299911678901748)
    ,(1500, "falling", -0.7388706207275391, 0.7278079986572266, 0.0533058643341064, 0.046578474342823, -0.1310305893421173, -0.0218384321779012)
    ,(1600, "falling", -0.7223987579345703, 0.7231974601745605, 0.1768441945314407, 0.0740674127425465, -0.09423323614256736, -0.02838341776458984)
    ,(1700, "falling", -0.8466672897338867, 1.015576839447022, 0.3397563397884369, 0.05015640013984282, -0.1209949524629682, -0.04998914045946923)
    ,(1800, "falling", -0.3906707763671875, 0.6379885673522949, 0.3954523801803589, 0.07590000810367709, -0.03761184641293151, 0.0287979330335341)
    ,(1900, "falling", -0.8676061630249023, 0.545799732208252, 0.2887592911720276, 0.1049888261726923, -0.01085013018122738, -0.02878338896802463)
    ,(2000, "falling", -0.6799111366271973, 0.3862724304199219, -0.3476599752902984, 0.1882119462603614, -0.2839359787425824, 0.3179844381908575)
    ,(2100, "falling", -0.1029138565063476, -0.1001324653625488, -1.190993785858154, -0.160526663064955, -0.3789109843117849, 0.2785327349390309)
    ,(2200, "falling", -0.4145522117614746, -0.0007905960083007, -1.845284700393677, -0.5090543727080018, -0.3627521395683289, 0.3927354415257764)
    ,(2300, "falling", -1.701473712921143, -0.7780666351318359, -1.666865825653076, -0.7832455890519271, -0.4365941172554375, 0.5598579928988499)
    ,(2400, "falling", -2.315230846405029, -2.501971244812012, -2.233747959136963, -1.023984696183884, 1.582228234836008, -0.6729553427014874)
    ,(2500, "falling", -2.985809087753296, -4.085277080535889, -1.856818199157715, -0.9998409918376387, 1.692903931651809, -2.806933039710625)
    ,(2600, "falling", -3.954225778579712, -3.951249122619629, -3.099813938140869, -0.4746423083401879, 0.7538440497148589, -2.680687563759945)
    ,(2700, "falling", -3.63595199584961, 2.219827651977539, 0.8921973705291748, -0.8093528002500497, 0.5344852464539692, -1.117512524127961)
    ,(2800, "falling", 1.556339263916016, -1.303864598274231, 0.8833200931549072, -1.273726730119617, 1.120043240842359, -0.5577272261892057)
    ,(2900, "falling", -5.200931549072266, -2.096135139465332, -4.002499103546143, -0.2288199303050852, 2.496860782305416, -0.04708753029505666)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
